{"js": "// Replace each \"NNN\u00d7N=\" multiplication expression in the document's\n// table cells with its updated value, per the commit diff. Every\n// expression in this worksheet is unique, so a body-wide search for the\n// exact old text unambiguously finds the single cell to update.\nconst replacements = [\n  [\"485\u00d78=\", \"359\u00d73=\"],\n  [\"345\u00d77=\", \"239\u00d77=\"],\n  [\"975\u00d77=\", \"647\u00d72=\"],\n  [\"550\u00d74=\", \"716\u00d76=\"],\n  [\"180\u00d79=\", \"299\u00d73=\"],\n  [\"754\u00d75=\", \"338\u00d73=\"],\n  [\"307\u00d76=\", \"808\u00d72=\"],\n  [\"666\u00d78=\", \"731\u00d78=\"],\n  [\"891\u00d75=\", \"757\u00d73=\"],\n  [\"754\u00d78=\", \"962\u00d76=\"],\n  [\"884\u00d73=\", \"813\u00d78=\"],\n  [\"965\u00d78=\", \"401\u00d74=\"],\n  [\"587\u00d74=\", \"822\u00d74=\"],\n  [\"819\u00d72=\", \"239\u00d76=\"],\n  [\"799\u00d73=\", \"720\u00d78=\"],\n  [\"171\u00d76=\", \"833\u00d79=\"],\n  [\"342\u00d78=\", \"762\u00d72=\"],\n  [\"842\u00d73=\", \"450\u00d77=\"],\n  [\"933\u00d76=\", \"473\u00d77=\"],\n  [\"223\u00d72=\", \"852\u00d74=\"],\n  [\"967\u00d77=\", \"480\u00d77=\"],\n  [\"746\u00d75=\", \"766\u00d73=\"],\n  [\"920\u00d75=\", \"733\u00d77=\"],\n  [\"703\u00d76=\", \"586\u00d72=\"],\n  [\"277\u00d79=\", \"885\u00d76=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each \"NNN\u00d7N=\" multiplication expression in the document's\n# table cells with its updated value, per the commit diff. Every\n# expression in this worksheet is unique, so a Find/Replace on the\n# exact old text unambiguously updates the single matching cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"485\u00d78=\", \"359\u00d73=\"),\n    @(\"345\u00d77=\", \"239\u00d77=\"),\n    @(\"975\u00d77=\", \"647\u00d72=\"),\n    @(\"550\u00d74=\", \"716\u00d76=\"),\n    @(\"180\u00d79=\", \"299\u00d73=\"),\n    @(\"754\u00d75=\", \"338\u00d73=\"),\n    @(\"307\u00d76=\", \"808\u00d72=\"),\n    @(\"666\u00d78=\", \"731\u00d78=\"),\n    @(\"891\u00d75=\", \"757\u00d73=\"),\n    @(\"754\u00d78=\", \"962\u00d76=\"),\n    @(\"884\u00d73=\", \"813\u00d78=\"),\n    @(\"965\u00d78=\", \"401\u00d74=\"),\n    @(\"587\u00d74=\", \"822\u00d74=\"),\n    @(\"819\u00d72=\", \"239\u00d76=\"),\n    @(\"799\u00d73=\", \"720\u00d78=\"),\n    @(\"171\u00d76=\", \"833\u00d79=\"),\n    @(\"342\u00d78=\", \"762\u00d72=\"),\n    @(\"842\u00d73=\", \"450\u00d77=\"),\n    @(\"933\u00d76=\", \"473\u00d77=\"),\n    @(\"223\u00d72=\", \"852\u00d74=\"),\n    @(\"967\u00d77=\", \"480\u00d77=\"),\n    @(\"746\u00d75=\", \"766\u00d73=\"),\n    @(\"920\u00d75=\", \"733\u00d77=\"),\n    @(\"703\u00d76=\", \"586\u00d72=\"),\n    @(\"277\u00d79=\", \"885\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n"}
